{"js": "// Insert a new paragraph right after the paragraph that ends with\n// \"...c\u00e1c trang web \u0111\u1ecba \u0111i\u1ec3m.\" and right before the paragraph that\n// begins with \"Trang web \u0111\u1ed3ng th\u1eddi c\u0169ng...\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"kh\u00f4ng c\u1ea7n b\u1ea1n m\u1ea5t th\u1eddi gian t\u00ecm ki\u1ebfm tr\u00ean c\u00e1c search engine hay c\u00e1c trang web \u0111\u1ecba \u0111i\u1ec3m.\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorText) !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Anchor paragraph not found\");\n}\n\nconst newText = \"Ngo\u00e0i ra, b\u1ea1n c\u00f3 th\u1ec3 k\u1ebft n\u1ed1i v\u1edbi m\u1ed9t ng\u01b0\u1eddi d\u00f9ng kh\u00e1c \u0111\u1ec3 ch\u01a1i caro ngay tr\u00ean trang web.\";\n// insertParagraph(\"After\") clones the anchor paragraph's mark formatting\n// (sz=24, szCs=24, lang=en-US), matching the rest of this section.\nanchor.insertParagraph(newText, \"After\");\n\nawait context.sync();\n", "ps1": "# Insert a new paragraph right after the paragraph that ends with\n# \"...c\u00e1c trang web \u0111\u1ecba \u0111i\u1ec3m.\" and right before the paragraph that\n# begins with \"Trang web \u0111\u1ed3ng th\u1eddi c\u0169ng...\".\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"c\u00e1c trang web \u0111\u1ecba \u0111i\u1ec3m.\")\nif (-not $found) {\n    throw \"Anchor text not found\"\n}\n\n$anchorPara = $rng.Paragraphs(1)\n\n$endRange = $anchorPara.Range\n$endRange.Collapse(0)            # wdCollapseEnd\n$endRange.InsertParagraphAfter()\n\n$newPara = $anchorPara.Next()\n$newPara.Range.Text = \"Ngo\u00e0i ra, b\u1ea1n c\u00f3 th\u1ec3 k\u1ebft n\u1ed1i v\u1edbi m\u1ed9t ng\u01b0\u1eddi d\u00f9ng kh\u00e1c \u0111\u1ec3 ch\u01a1i caro ngay tr\u00ean trang web.\"\n"}
